# Weekly update: a new Coliflor price record (week of 2021-10-08) is added
# at the top of the data block (row 53), pushing every existing record
# down by one row. The former last record (old row 148) ends up as the
# new last record (row 149).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 53; Excel shifts rows 53..148 down to 54..149
# and also extends the used range / dimension automatically, copying down
# the number format (date style) already present on column D.
$ws.Rows(53).Insert()

# Populate the new row 53 with the new weekly record.
$ws.Range("A53").Value2 = 5
$ws.Range("B53").Value2 = "Macroferia Regional de Talca"
$ws.Range("C53").Value2 = "Maule"
$ws.Range("D53").Value2 = 44477
$ws.Range("E53").Value2 = 7
$ws.Range("F53").Value2 = 100112008
$ws.Range("G53").Value2 = "Coliflor"
$ws.Range("H53").Value2 = "Sin especificar"
$ws.Range("I53").Value2 = "Primera"
$ws.Range("J53").Value2 = 3000
$ws.Range("K53").Value2 = 600
$ws.Range("L53").Value2 = 600
$ws.Range("M53").Value2 = 600
$ws.Range("N53").Value2 = "`$/unidad"
$ws.Range("O53").Value2 = "Región del Maule"
$ws.Range("P53").Value2 = 600
$ws.Range("Q53").Value2 = 1
$ws.Range("R53").Value2 = "Hortaliza"
